$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 107: date and value were updated
$ws.Cells.Item(107, 1).Value = 45747
$ws.Cells.Item(107, 2).Value = 1.85

# Add new row 108 with the same formatting as the prior rows
$ws.Cells.Item(108, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(108, 2).NumberFormat = "0.00"
$ws.Cells.Item(108, 1).Value = 45777
$ws.Cells.Item(108, 2).Value = 1.94

# Update the active selection to match the new state (B109, one row below the new data)
$ws.Range("B109").Select()
